# Update metricas_recorrencia_anual data for the 2025 row (row 6):
#  - total_customers (C6): 405 -> 406
#  - new_customers   (E6): 98  -> 99
#  - new_rate        (G6): recalculated as new_customers / total_customers * 100
#  - returning_rate  (H6): recalculated as 100 - new_rate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 406
$ws.Range("E6").Value = 99
$ws.Range("G6").Value = 24.38423645320197
$ws.Range("H6").Value = 75.61576354679804
